# Update the "Förändrad" (Changed) date column (C) for all data rows
# (rows 2-101) from 2023-09-15 (45184) to 2023-09-16 (45185).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C101").Value = 45185
